$d = $word.ActiveDocument

# 1. Update the title text.
[void]$d.Content.Find.Execute("Complex Test Document", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Test Document with Table", 2)

# 2. Update the intro paragraph text.
[void]$d.Content.Find.Execute("This document has multiple tables.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This is a test document.", 2)

# 3. Remove the first (Key/Value, Status/Active) table entirely.
$d.Tables.Item(1).Delete()

# 4. Remove the now-orphaned "Here is another table:" paragraph (text + mark).
$rng = $d.Content
[void]$rng.Find.Execute("Here is another table:", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$d.Range($rng.Start, $rng.End + 1).Delete()

# 5. Apply the "Light Grid Accent 1" table style to the remaining table.
$t = $d.Tables.Item(1)
$t.Style = "Light Grid Accent 1"

# 6. Rewrite the header row and data cells in place.
$t.Cell(1, 1).Range.Text = "Name"
$t.Cell(1, 2).Range.Text = "Age"
$t.Cell(1, 3).Range.Text = "City"

$t.Cell(2, 1).Range.Text = "Alice"
$t.Cell(2, 2).Range.Text = "30"
$t.Cell(2, 3).Range.Text = "NYC"

$t.Cell(3, 1).Range.Text = "Bob"
$t.Cell(3, 2).Range.Text = "25"
$t.Cell(3, 3).Range.Text = "LA"

# 7. Drop the trailing "Orange" row.
$t.Rows.Item(4).Delete()
